$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells (styled like the existing header row)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the neighboring header cell (AC1) so the new
# headers match the rest of the header row (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($row = 2; $row -le 63; $row++) {
    $ws.Cells.Item($row, 30).Value = 67
    $ws.Cells.Item($row, 31).Value = 95
    $ws.Cells.Item($row, 32).Value = 0
}
